$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new employee row (Hank Pym) below the existing data
$ws.Range("A15").Value = 451
$ws.Range("B15").Value = "Hank"
$ws.Range("C15").Value = "Pym"
$ws.Range("D15").Value = 35000
$ws.Range("E15").Value = 311

# Update the active selection to reflect where the user ended up after entry
$ws.Range("E16").Select()
